$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '59.036.19'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -3.13%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '3.236.25'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -3.61%  '
$ws.Range('E4').Value = '  +0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '538.74'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -5.35%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '136.71'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -9.18%  '
$ws.Range('E7').Value = '  -0.11%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '3.235.37'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -3.59%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.460'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -4.29%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '7.63'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -4.61%  '
$ws.Range('E11').Value = '  -5.63%  '
$ws.Range('E12').Value = '  -4.76%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '3.790.12'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -3.47%  '
$ws.Range('E14').Value = '  -1.36%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '26.20'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -7.03%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '3.233.75'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -3.68%  '
$ws.Range('E17').Value = '  -6.30%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '59.084.94'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -3.20%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '5.90'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -6.89%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '13.23'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -6.78%  '
$ws.Range('E21').Value = '  -6.32%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '361.69'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.01%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '70.60'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -6.28%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.521'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -7.22%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '3.369.72'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -4.20%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '0.0₃0976'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -10.42%  '
$ws.Range('E28').Value = '  -3.54%  '
$ws.Range('E29').Value = '  +0.03%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '7.03'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -5.12%  '
$ws.Range('E31').Value = '  -0.06%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.93'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -7.14%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '7.04'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -8.60%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '22.04'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -3.76%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '163.81'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -3.45%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '4.94'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -8.49%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '6.40'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -5.49%  '
$ws.Range('E39').Value = '  -6.66%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '26.14'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -10.73%  '
$ws.Range('E41').Value = '  -6.25%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '3.266.33'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -3.77%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '41.17'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('E44').Value = '  -5.53%  '
$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '1.10'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -3.91%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '4.02'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -6.69%  '
$ws.Range('E47').Value = '  -6.48%  '
$ws.Range('E48').Value = '  +0.01%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '2.296.58'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -8.66%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '6.29'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -6.41%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '20.73'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -9.37%  '
